$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("G2").Value = 42.32036466666667
$ws.Range("H2").Value = 126.961094
$ws.Range("I2").Value = 0.285778576657872
$ws.Range("J2").Value = 0.2880046678857171
$ws.Range("M2").Value = 6.712486666666666
$ws.Range("N2").Value = 20.13746
$ws.Range("O2").Value = 0.6330487633990675
$ws.Range("P2").Value = 0.6414503882251803
$ws.Range("Q2").Value = 284.0748835534711
$ws.Range("R2").Value = 2556.673951981239
$ws.Range("S2").Value = 0.1809117745592115
$ws.Range("T2").Value = 0.1847407060259574
$ws.Range("G3").Value = 42.32036466666667
$ws.Range("H3").Value = 126.961094
$ws.Range("I3").Value = 0.285778576657872
$ws.Range("J3").Value = 0.2880046678857171
$ws.Range("O3").Value = 0.290741083484562
$ws.Range("P3").Value = 0.2945997080427384
$ws.Range("Q3").Value = 130.4674208533627
$ws.Range("R3").Value = 1174.206787680264
$ws.Range("S3").Value = 0.08308757301418566
$ws.Range("T3").Value = 0.08484609107407809
$ws.Range("G4").Value = 42.32036466666667
$ws.Range("H4").Value = 126.961094
$ws.Range("I4").Value = 0.285778576657872
$ws.Range("J4").Value = 0.2880046678857171
$ws.Range("M4").Value = 0.2495096666666667
$ws.Range("N4").Value = 0.748529
$ws.Range("O4").Value = 0.02353103905946135
$ws.Range("P4").Value = 0.02384333563656022
$ws.Range("Q4").Value = 10.55934008119178
$ws.Range("R4").Value = 95.034060730726
$ws.Range("S4").Value = 0.006724666849693657
$ws.Range("T4").Value = 0.00686699196129521
$ws.Range("G5").Value = 42.32036466666667
$ws.Range("H5").Value = 126.961094
$ws.Range("I5").Value = 0.285778576657872
$ws.Range("J5").Value = 0.2880046678857171
$ws.Range("M5").Value = 0.4166465
$ws.Range("N5").Value = 0.8332930000000001
$ws.Range("O5").Value = 0.03929356804674715
$ws.Range("P5").Value = 0.02654337331298611
$ws.Range("Q5").Value = 17.63263181709033
$ws.Range("R5").Value = 105.795790902542
$ws.Range("S5").Value = 0.01122925994820864
$ws.Range("T5").Value = 0.007644615415573172
$ws.Range("G6").Value = 42.32036466666667
$ws.Range("H6").Value = 126.961094
$ws.Range("I6").Value = 0.285778576657872
$ws.Range("J6").Value = 0.2880046678857171
$ws.Range("K6").Value = 2
$ws.Range("L6").Value = 0.6666666666666666
$ws.Range("M6").Value = 0.1419326666666667
$ws.Range("N6").Value = 0.425798
$ws.Range("O6").Value = 0.01338554601016197
$ws.Range("P6").Value = 0.01356319478253491
$ws.Range("Q6").Value = 6.006642211445778
$ws.Range("R6").Value = 54.059779903012
$ws.Range("S6").Value = 0.003825302286572545
$ws.Range("T6").Value = 0.003906263408813257
$ws.Range("I7").Value = 0.04213668412459876
$ws.Range("J7").Value = 0.04246491062777905
$ws.Range("M7").Value = 6.712486666666666
$ws.Range("N7").Value = 20.13746
$ws.Range("O7").Value = 0.6330487633990675
$ws.Range("P7").Value = 0.6414503882251803
$ws.Range("Q7").Value = 41.88548272586221
$ws.Range("R7").Value = 376.9693445327599
$ws.Range("S7").Value = 0.02667457577881436
$ws.Range("T7").Value = 0.02723913340813645
$ws.Range("I8").Value = 0.04213668412459876
$ws.Range("J8").Value = 0.04246491062777905
$ws.Range("O8").Value = 0.290741083484562
$ws.Range("P8").Value = 0.2945997080427384
$ws.Range("S8").Value = 0.01225086519683258
$ws.Range("T8").Value = 0.01251015027300469
$ws.Range("I9").Value = 0.04213668412459876
$ws.Range("J9").Value = 0.04246491062777905
$ws.Range("M9").Value = 0.2495096666666667
$ws.Range("N9").Value = 0.748529
$ws.Range("O9").Value = 0.02353103905946135
$ws.Range("P9").Value = 0.02384333563656022
$ws.Range("Q9").Value = 1.556924185041555
$ws.Range("R9").Value = 14.012317665374
$ws.Range("S9").Value = 0.0009915199599721187
$ws.Range("T9").Value = 0.001012505116874669
$ws.Range("I10").Value = 0.04213668412459876
$ws.Range("J10").Value = 0.04246491062777905
$ws.Range("M10").Value = 0.4166465
$ws.Range("N10").Value = 0.8332930000000001
$ws.Range("O10").Value = 0.03929356804674715
$ws.Range("P10").Value = 0.02654337331298611
$ws.Range("Q10").Value = 2.599847216859667
$ws.Range("R10").Value = 15.599083301158
$ws.Range("S10").Value = 0.001655700664914212
$ws.Range("T10").Value = 0.001127161975495731
$ws.Range("I11").Value = 0.04213668412459876
$ws.Range("J11").Value = 0.04246491062777905
$ws.Range("K11").Value = 2
$ws.Range("L11").Value = 0.6666666666666666
$ws.Range("M11").Value = 0.1419326666666667
$ws.Range("N11").Value = 0.425798
$ws.Range("O11").Value = 0.01338554601016197
$ws.Range("P11").Value = 0.01356319478253491
$ws.Range("Q11").Value = 0.8856506616875555
$ws.Range("R11").Value = 7.970855955187999
$ws.Range("S11").Value = 0.000564022524065478
$ws.Range("T11").Value = 0.0005759598542675038
$ws.Range("G12").Value = 42.241047
$ws.Range("H12").Value = 126.723141
$ws.Range("I12").Value = 0.2852429647825406
$ws.Range("J12").Value = 0.2874648838260633
$ws.Range("M12").Value = 6.712486666666666
$ws.Range("N12").Value = 20.13746
$ws.Range("O12").Value = 0.6330487633990675
$ws.Range("P12").Value = 0.6414503882251803
$ws.Range("Q12").Value = 283.54246477354
$ws.Range("R12").Value = 2551.88218296186
$ws.Range("S12").Value = 0.1805727061238711
$ws.Range("T12").Value = 0.1843944613313346
$ws.Range("G13").Value = 42.241047
$ws.Range("H13").Value = 126.723141
$ws.Range("I13").Value = 0.2852429647825406
$ws.Range("J13").Value = 0.2874648838260633
$ws.Range("O13").Value = 0.290741083484562
$ws.Range("P13").Value = 0.2945997080427384
$ws.Range("Q13").Value = 130.222896226044
$ws.Range("R13").Value = 1172.006066034396
$ws.Range("S13").Value = 0.08293184863722459
$ws.Range("T13").Value = 0.08468707084769794
$ws.Range("G14").Value = 42.241047
$ws.Range("H14").Value = 126.723141
$ws.Range("I14").Value = 0.2852429647825406
$ws.Range("J14").Value = 0.2874648838260633
$ws.Range("M14").Value = 0.2495096666666667
$ws.Range("N14").Value = 0.748529
$ws.Range("O14").Value = 0.02353103905946135
$ws.Range("P14").Value = 0.02384333563656022
$ws.Range("Q14").Value = 10.539549556621
$ws.Range("R14").Value = 94.855946009589
$ws.Range("S14").Value = 0.006712063345734522
$ws.Range("T14").Value = 0.006854121708789619
$ws.Range("G15").Value = 42.241047
$ws.Range("H15").Value = 126.723141
$ws.Range("I15").Value = 0.2852429647825406
$ws.Range("J15").Value = 0.2874648838260633
$ws.Range("M15").Value = 0.4166465
$ws.Range("N15").Value = 0.8332930000000001
$ws.Range("O15").Value = 0.03929356804674715
$ws.Range("P15").Value = 0.02654337331298611
$ws.Range("Q15").Value = 17.5995843888855
$ws.Range("R15").Value = 105.597506333313
$ws.Range("S15").Value = 0.01120821384653866
$ws.Range("T15").Value = 0.007630287725769381
$ws.Range("G16").Value = 42.241047
$ws.Range("H16").Value = 126.723141
$ws.Range("I16").Value = 0.2852429647825406
$ws.Range("J16").Value = 0.2874648838260633
$ws.Range("K16").Value = 2
$ws.Range("L16").Value = 0.6666666666666666
$ws.Range("M16").Value = 0.1419326666666667
$ws.Range("N16").Value = 0.425798
$ws.Range("O16").Value = 0.01338554601016197
$ws.Range("P16").Value = 0.01356319478253491
$ws.Range("Q16").Value = 5.995384443502001
$ws.Range("R16").Value = 53.958459991518
$ws.Range("S16").Value = 0.003818132829171706
$ws.Range("T16").Value = 0.003898942212471664
$ws.Range("G17").Value = 3.4338745
$ws.Range("H17").Value = 6.867749
$ws.Range("I17").Value = 0.02318807445921414
$ws.Range("J17").Value = 0.0155791330048516
$ws.Range("M17").Value = 6.712486666666666
$ws.Range("N17").Value = 20.13746
$ws.Range("O17").Value = 0.6330487633990675
$ws.Range("P17").Value = 0.6414503882251803
$ws.Range("Q17").Value = 23.04983679625666
$ws.Range("R17").Value = 138.29902077754
$ws.Range("S17").Value = 0.01467918186201101
$ws.Range("T17").Value = 0.009993240914173775
$ws.Range("G18").Value = 3.4338745
$ws.Range("H18").Value = 6.867749
$ws.Range("I18").Value = 0.02318807445921414
$ws.Range("J18").Value = 0.0155791330048516
$ws.Range("O18").Value = 0.290741083484562
$ws.Range("P18").Value = 0.2945997080427384
$ws.Range("Q18").Value = 10.586126870074
$ws.Range("R18").Value = 63.516761220444
$ws.Range("S18").Value = 0.006741725892192616
$ws.Range("T18").Value = 0.00458960803478827
$ws.Range("G19").Value = 3.4338745
$ws.Range("H19").Value = 6.867749
$ws.Range("I19").Value = 0.02318807445921414
$ws.Range("J19").Value = 0.0155791330048516
$ws.Range("M19").Value = 0.2495096666666667
$ws.Range("N19").Value = 0.748529
$ws.Range("O19").Value = 0.02353103905946135
$ws.Range("P19").Value = 0.02384333563656022
$ws.Range("Q19").Value = 0.8567848818701667
$ws.Range("R19").Value = 5.140709291221
$ws.Range("S19").Value = 0.0005456394858134661
$ws.Range("T19").Value = 0.0003714584971612896
$ws.Range("G20").Value = 3.4338745
$ws.Range("H20").Value = 6.867749
$ws.Range("I20").Value = 0.02318807445921414
$ws.Range("J20").Value = 0.0155791330048516
$ws.Range("M20").Value = 0.4166465
$ws.Range("N20").Value = 0.8332930000000001
$ws.Range("O20").Value = 0.03929356804674715
$ws.Range("P20").Value = 0.02654337331298611
$ws.Range("Q20").Value = 1.43071179186425
$ws.Range("R20").Value = 5.722847167457
$ws.Range("S20").Value = 0.0009111421816361704
$ws.Range("T20").Value = 0.000413522743240439
$ws.Range("G21").Value = 3.4338745
$ws.Range("H21").Value = 6.867749
$ws.Range("I21").Value = 0.02318807445921414
$ws.Range("J21").Value = 0.0155791330048516
$ws.Range("K21").Value = 2
$ws.Range("L21").Value = 0.6666666666666666
$ws.Range("M21").Value = 0.1419326666666667
$ws.Range("N21").Value = 0.425798
$ws.Range("O21").Value = 0.01338554601016197
$ws.Range("P21").Value = 0.01356319478253491
$ws.Range("Q21").Value = 0.4873789647836667
$ws.Range("R21").Value = 2.924273788702
$ws.Range("S21").Value = 0.0003103850375608724
$ws.Range("T21").Value = 0.0002113028154878205
$ws.Range("G22").Value = 53.85273233333334
$ws.Range("H22").Value = 161.558197
$ws.Range("I22").Value = 0.3636536999757743
$ws.Range("J22").Value = 0.3664864046555889
$ws.Range("M22").Value = 6.712486666666666
$ws.Range("N22").Value = 20.13746
$ws.Range("O22").Value = 0.6330487633990675
$ws.Range("P22").Value = 0.6414503882251803
$ws.Range("Q22").Value = 361.4857477510689
$ws.Range("R22").Value = 3253.371729759619
$ws.Range("S22").Value = 0.2302105250751594
$ws.Range("T22").Value = 0.235082846545578
$ws.Range("G23").Value = 53.85273233333334
$ws.Range("H23").Value = 161.558197
$ws.Range("I23").Value = 0.3636536999757743
$ws.Range("J23").Value = 0.3664864046555889
$ws.Range("O23").Value = 0.290741083484562
$ws.Range("P23").Value = 0.2945997080427384
$ws.Range("Q23").Value = 166.0200035792814
$ws.Range("R23").Value = 1494.180032213532
$ws.Range("S23").Value = 0.1057290707441264
$ws.Range("T23").Value = 0.1079667878131694
$ws.Range("G24").Value = 53.85273233333334
$ws.Range("H24").Value = 161.558197
$ws.Range("I24").Value = 0.3636536999757743
$ws.Range("J24").Value = 0.3664864046555889
$ws.Range("M24").Value = 0.2495096666666667
$ws.Range("N24").Value = 0.748529
$ws.Range("O24").Value = 0.02353103905946135
$ws.Range("P24").Value = 0.02384333563656022
$ws.Range("Q24").Value = 13.43677729357922
$ws.Range("R24").Value = 120.930995642213
$ws.Range("S24").Value = 0.008557149418247586
$ws.Range("T24").Value = 0.008738258352439434
$ws.Range("G25").Value = 53.85273233333334
$ws.Range("H25").Value = 161.558197
$ws.Range("I25").Value = 0.3636536999757743
$ws.Range("J25").Value = 0.3664864046555889
$ws.Range("M25").Value = 0.4166465
$ws.Range("N25").Value = 0.8332930000000001
$ws.Range("O25").Value = 0.03929356804674715
$ws.Range("P25").Value = 0.02654337331298611
$ws.Range("Q25").Value = 22.43755244212017
$ws.Range("R25").Value = 134.625314652721
$ws.Range("S25").Value = 0.01428925140544946
$ws.Range("T25").Value = 0.009727785452907388
$ws.Range("G26").Value = 53.85273233333334
$ws.Range("H26").Value = 161.558197
$ws.Range("I26").Value = 0.3636536999757743
$ws.Range("J26").Value = 0.3664864046555889
$ws.Range("K26").Value = 2
$ws.Range("L26").Value = 0.6666666666666666
$ws.Range("M26").Value = 0.1419326666666667
$ws.Range("N26").Value = 0.425798
$ws.Range("O26").Value = 0.01338554601016197
$ws.Range("P26").Value = 0.01356319478253491
$ws.Range("Q26").Value = 7.643461907356223
$ws.Range("R26").Value = 68.79115716620601
$ws.Range("S26").Value = 0.004867703332791363
$ws.Range("T26").Value = 0.004970726491494659
